$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 649.4167
$ws.Range("I28").Value = 479.3
$ws.Range("K28").Value = 479.3
$ws.Range("M28").Value = 5.699999999999989

$ws.Range("H33").Value = 2820.2
$ws.Range("I33").Value = 799.25
$ws.Range("J33").Value = 4167.5
$ws.Range("K33").Value = 799.25
$ws.Range("L33").Value = 4167.5
$ws.Range("M33").Value = -570.25
$ws.Range("N33").Value = -4625.5

$ws.Range("H38").Value = 1198.1428
$ws.Range("I38").Value = 64.666664
$ws.Range("K38").Value = 193.999992
$ws.Range("M38").Value = 178.000008

$ws.Range("H39").Value = 30.875
$ws.Range("I39").Value = 30.875
$ws.Range("K39").Value = 92.625
$ws.Range("M39").Value = 203.375

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H70").Value = 3230.5

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H73").Value = 3230.5

$ws.Range("H86").Value = 11877.667
$ws.Range("J86").Value = 5800
$ws.Range("L86").Value = 5800
$ws.Range("N86").Value = -8046

$ws.Range("H89").Value = 11877.667
$ws.Range("J89").Value = 5800
$ws.Range("L89").Value = 29000
$ws.Range("N89").Value = -40232

$ws.Range("H111").Value = 1118.1428
$ws.Range("I111").Value = 1118.1428
$ws.Range("K111").Value = 3354.4284
$ws.Range("M111").Value = -287.4284000000002

$ws.Range("H112").Value = 2475
$ws.Range("I112").Value = 1700
$ws.Range("K112").Value = 5100
$ws.Range("M112").Value = -3992

$ws.Range("H116").Value = 6466.6
$ws.Range("I116").Value = 5800
$ws.Range("K116").Value = 5800
$ws.Range("M116").Value = -2358

$ws.Range("H137").Value = 2470.2666
$ws.Range("I137").Value = 1694.625
$ws.Range("J137").Value = 3356.7144
$ws.Range("K137").Value = 5083.875
$ws.Range("L137").Value = 10070.1432
$ws.Range("M137").Value = -2533.875
$ws.Range("N137").Value = -15170.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6432.2925
$ws.Range("I32").Value = 5326.2563
$ws.Range("K32").Value = 5326.2563
$ws.Range("M32").Value = -5039.2563

$ws.Range("H63").Value = 2515.2856
$ws.Range("I63").Value = 1402
$ws.Range("J63").Value = 3999.6667
$ws.Range("K63").Value = 1402
$ws.Range("L63").Value = 3999.6667
$ws.Range("M63").Value = -716
$ws.Range("N63").Value = -5371.6667

$ws.Range("H66").Value = 2515.2856
$ws.Range("I66").Value = 1402
$ws.Range("J66").Value = 3999.6667
$ws.Range("K66").Value = 7010
$ws.Range("L66").Value = 19998.3335
$ws.Range("M66").Value = -3578
$ws.Range("N66").Value = -26862.3335

$ws.Range("H110").Value = 2351.55
$ws.Range("I110").Value = 631.1818
$ws.Range("K110").Value = 631.1818
$ws.Range("M110").Value = 1413.8182

$ws.Range("H122").Value = 2129.5
$ws.Range("I122").Value = 2129.5
$ws.Range("K122").Value = 6388.5
$ws.Range("M122").Value = -3938.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 766.6667
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H107").Value = 874.2857
$ws.Range("I107").Value = 813.3333
$ws.Range("J107").Value = 920
$ws.Range("K107").Value = 813.3333
$ws.Range("L107").Value = 920
$ws.Range("M107").Value = 1106.6667
$ws.Range("N107").Value = -4760

$ws.Range("H120").Value = 19999.25
$ws.Range("J120").Value = 20000
$ws.Range("L120").Value = 20000
$ws.Range("N120").Value = -27258

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 752.1111
$ws.Range("I5").Value = 469.42856
$ws.Range("J5").Value = 1741.5
$ws.Range("K5").Value = 1408.28568
$ws.Range("L5").Value = 5224.5
$ws.Range("M5").Value = -1296.28568
$ws.Range("N5").Value = -5448.5

$ws.Range("H135").Value = 752.1111
$ws.Range("I135").Value = 469.42856
$ws.Range("J135").Value = 1741.5
$ws.Range("K135").Value = 4224.85704
$ws.Range("L135").Value = 15673.5
$ws.Range("M135").Value = -1689.85704
$ws.Range("N135").Value = -20743.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 8164.75
$ws.Range("I99").Value = 8164.75
$ws.Range("K99").Value = 8164.75
$ws.Range("M99").Value = -5918.75

$ws.Range("H102").Value = 1478.6364
$ws.Range("I102").Value = 1140.5555
$ws.Range("K102").Value = 1140.5555
$ws.Range("M102").Value = 481.4445000000001

$ws.Range("H107").Value = 557.6
$ws.Range("I107").Value = 557.6
$ws.Range("K107").Value = 557.6
$ws.Range("M107").Value = 1362.4

$ws.Range("H113").Value = 1534.2222
$ws.Range("I113").Value = 1534.2222
$ws.Range("K113").Value = 1534.2222
$ws.Range("M113").Value = 635.7778000000001

$ws.Range("H122").Value = 17893430
$ws.Range("I122").Value = 20858834
$ws.Range("K122").Value = 62576502
$ws.Range("M122").Value = -62574052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 434
$ws.Range("I9").Value = 553.6667
$ws.Range("J9").Value = 75
$ws.Range("K9").Value = 553.6667
$ws.Range("L9").Value = 75
$ws.Range("M9").Value = -329.6667
$ws.Range("N9").Value = -523

$ws.Range("H40").Value = 6306.2144
$ws.Range("I40").Value = 6123.9165
$ws.Range("K40").Value = 6123.9165
$ws.Range("M40").Value = -5987.9165

$ws.Range("H46").Value = 3238.2307
$ws.Range("I46").Value = 799.6667
$ws.Range("J46").Value = 3969.8
$ws.Range("K46").Value = 799.6667
$ws.Range("L46").Value = 3969.8
$ws.Range("M46").Value = -611.6667
$ws.Range("N46").Value = -4345.8

$ws.Range("H122").Value = 3529.818
$ws.Range("I122").Value = 3094.25
$ws.Range("J122").Value = 4691.3335
$ws.Range("K122").Value = 9282.75
$ws.Range("L122").Value = 14074.0005
$ws.Range("M122").Value = -6832.75
$ws.Range("N122").Value = -18974.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 473
$ws.Range("I81").Value = 473
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 946
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 115
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 473
$ws.Range("I84").Value = 473
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4730
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 574
$ws.Range("N84").ClearContents()

$ws.Range("H132").Value = 1044.2222
$ws.Range("I132").Value = 1044.2222
$ws.Range("K132").Value = 3132.6666
$ws.Range("M132").Value = -602.6665999999996
